$d = $word.ActiveDocument

# "Nama Lembaga" row, value cell currently reads ": ${nama_lembaga}" with
# the leading ": " and the placeholder-opener "${" split across two
# separate (but identically-formatted) runs. Merge just those first four
# characters into a single run so it reads "${" in the same run as ": ",
# matching how the other fields in this table are authored, while leaving
# the "nama_lembaga" / "}" runs (and their spell-check proofErr markers)
# completely untouched.
$t = $d.Tables.Item(3)
$cell = $t.Cell(2, 2)
$rng = $cell.Range

$sub = $d.Range($rng.Start, $rng.Start + 4)
if ($sub.Text -ne ": $" + "{") {
    throw "unexpected cell prefix: [" + $sub.Text + "]"
}

# Re-set through an intermediate value so the host actually re-splits the
# run boundaries instead of treating the assignment as a no-op (the text
# content before/after is identical, only the run split changes).
$sub.Text = "~~~~"
$sub2 = $d.Range($rng.Start, $rng.Start + 4)
$sub2.Text = ": $" + "{"
